$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 37323.875
$ws.Range("J3").Value = 37323.875
$ws.Range("L3").Value = 37323.875
$ws.Range("N3").Value = -37551.875
$ws.Range("H6").Value = 1825
$ws.Range("I6").Value = 1825
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 5475
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -5363
$ws.Range("N6").ClearContents()
$ws.Range("H102").Value = 37323.875
$ws.Range("J102").Value = 37323.875
$ws.Range("L102").Value = 37323.875
$ws.Range("N102").Value = -43813.875
$ws.Range("H113").Value = 3480.2666
$ws.Range("I113").Value = 3920.4
$ws.Range("J113").Value = 2600
$ws.Range("K113").Value = 3920.4
$ws.Range("L113").Value = 2600
$ws.Range("M113").Value = -666.4000000000001
$ws.Range("N113").Value = -9108
$ws.Range("H129").Value = 1515.2812
$ws.Range("I129").Value = 1747.3636
$ws.Range("K129").Value = 5242.0908
$ws.Range("M129").Value = -242.0907999999999
$ws.Range("H132").Value = 20808.291
$ws.Range("I132").Value = 2843.4634
$ws.Range("J132").Value = 126030.86
$ws.Range("K132").Value = 8530.3902
$ws.Range("L132").Value = 378092.58
$ws.Range("M132").Value = -6000.3902
$ws.Range("N132").Value = -383152.58
$ws.Range("H136").Value = 36049.617
$ws.Range("J136").Value = 36049.617
$ws.Range("L136").Value = 36049.617
$ws.Range("N136").Value = -46249.617
$ws.Range("H137").Value = 3543.2205
$ws.Range("I137").Value = 1184.1538
$ws.Range("K137").Value = 3552.4614
$ws.Range("M137").Value = -1002.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 5000
$ws.Range("J8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("N8").Value = -5288
$ws.Range("H32").Value = 24917.027
$ws.Range("I32").Value = 24346.463
$ws.Range("K32").Value = 24346.463
$ws.Range("M32").Value = -24059.463
$ws.Range("H132").Value = 16130677
$ws.Range("I132").Value = 21740238
$ws.Range("J132").Value = 3190
$ws.Range("K132").Value = 65220714
$ws.Range("L132").Value = 9570
$ws.Range("M132").Value = -65218184
$ws.Range("N132").Value = -14630

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 650
$ws.Range("J12").Value = 650
$ws.Range("L12").Value = 650
$ws.Range("N12").Value = -986
$ws.Range("H20").Value = 2595.6538
$ws.Range("I20").Value = 2163.75
$ws.Range("J20").Value = 3286.7
$ws.Range("K20").Value = 2163.75
$ws.Range("L20").Value = 3286.7
$ws.Range("M20").Value = -1916.75
$ws.Range("N20").Value = -3780.7
$ws.Range("H92").Value = 46130.668
$ws.Range("J92").Value = 46130.668
$ws.Range("L92").Value = 46130.668
$ws.Range("N92").Value = -51122.668
$ws.Range("H95").Value = 42163.25
$ws.Range("J95").Value = 42163.25
$ws.Range("L95").Value = 42163.25
$ws.Range("N95").Value = -47655.25
$ws.Range("H100").Value = 23814.334
$ws.Range("J100").Value = 23814.334
$ws.Range("L100").Value = 23814.334
$ws.Range("N100").Value = -25978.334
$ws.Range("H134").Value = 3748.0107
$ws.Range("I134").Value = 1641.9062
$ws.Range("J134").Value = 4852.8525
$ws.Range("K134").Value = 4925.7186
$ws.Range("L134").Value = 14558.5575
$ws.Range("M134").Value = -2390.7186
$ws.Range("N134").Value = -19628.5575

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 200
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H99").Value = 1681.6666
$ws.Range("I99").Value = 1322.2
$ws.Range("J99").Value = 1861.4
$ws.Range("K99").Value = 1322.2
$ws.Range("L99").Value = 1861.4
$ws.Range("M99").Value = 175.8
$ws.Range("N99").Value = -4857.4
$ws.Range("H126").Value = 1681.6666
$ws.Range("I126").Value = 1322.2
$ws.Range("J126").Value = 1861.4
$ws.Range("K126").Value = 3966.6
$ws.Range("L126").Value = 5584.200000000001
$ws.Range("M126").Value = -1496.6
$ws.Range("N126").Value = -10524.2
$ws.Range("H132").Value = 36135.78
$ws.Range("I132").Value = 1545.7667
$ws.Range("J132").Value = 130472.18
$ws.Range("K132").Value = 4637.300099999999
$ws.Range("L132").Value = 391416.54
$ws.Range("M132").Value = -2107.300099999999
$ws.Range("N132").Value = -396476.54
$ws.Range("H134").Value = 779148.3
$ws.Range("I134").Value = 1092.3
$ws.Range("J134").Value = 1751718.4
$ws.Range("K134").Value = 3276.9
$ws.Range("L134").Value = 5255155.199999999
$ws.Range("M134").Value = -741.8999999999996
$ws.Range("N134").Value = -5260225.199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 43987.53
$ws.Range("I131").Value = 9678.637000000001
$ws.Range("J131").Value = 53919.05
$ws.Range("K131").Value = 29035.911
$ws.Range("L131").Value = 161757.15
$ws.Range("M131").Value = -23995.911
$ws.Range("N131").Value = -171837.15
$ws.Range("H136").Value = 50002016
$ws.Range("I136").Value = 62501644
$ws.Range("K136").Value = 187504932
$ws.Range("M136").Value = -187499832

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4580
$ws.Range("I70").Value = 4659.2593
$ws.Range("J70").Value = 4223.3335
$ws.Range("K70").Value = 4659.2593
$ws.Range("L70").Value = 4223.3335
$ws.Range("M70").Value = -4389.2593
$ws.Range("N70").Value = -4763.3335
$ws.Range("H73").Value = 4580
$ws.Range("I73").Value = 4659.2593
$ws.Range("J73").Value = 4223.3335
$ws.Range("K73").Value = 4659.2593
$ws.Range("L73").Value = 4223.3335
$ws.Range("M73").Value = -3723.2593
$ws.Range("N73").Value = -6095.3335
$ws.Range("H95").Value = 30974.5
$ws.Range("J95").Value = 30974.5
$ws.Range("L95").Value = 30974.5
$ws.Range("N95").Value = -36466.5
$ws.Range("H98").Value = 38293.8
$ws.Range("J98").Value = 38293.8
$ws.Range("L98").Value = 38293.8
$ws.Range("N98").Value = -44283.8
$ws.Range("H132").Value = 3984.5588
$ws.Range("I132").Value = 1820.0588
$ws.Range("J132").Value = 6149.0586
$ws.Range("K132").Value = 5460.1764
$ws.Range("L132").Value = 18447.1758
$ws.Range("M132").Value = -2930.1764
$ws.Range("N132").Value = -23507.1758
$ws.Range("H136").Value = 19201.875
$ws.Range("J136").Value = 19201.875
$ws.Range("L136").Value = 57605.625
$ws.Range("N136").Value = -62705.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 10000
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 10000
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H104").Value = 24117.666
$ws.Range("J104").Value = 24117.666
$ws.Range("L104").Value = 24117.666
$ws.Range("N104").Value = -31105.666
$ws.Range("H106").Value = 44365.5
$ws.Range("J106").Value = 44365.5
$ws.Range("L106").Value = 44365.5
$ws.Range("N106").Value = -46889.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 90004
$ws.Range("J3").Value = 90004
$ws.Range("L3").Value = 90004
$ws.Range("N3").Value = -90232
$ws.Range("H69").Value = 19756.834
$ws.Range("J69").Value = 19756.834
$ws.Range("L69").Value = 19756.834
$ws.Range("N69").Value = -21254.834
$ws.Range("H72").Value = 19756.834
$ws.Range("J72").Value = 19756.834
$ws.Range("L72").Value = 59270.50199999999
$ws.Range("N72").Value = -66758.50199999999
$ws.Range("H80").Value = 35615
$ws.Range("J80").Value = 35615
$ws.Range("L80").Value = 35615
$ws.Range("N80").Value = -37611
$ws.Range("H83").Value = 35615
$ws.Range("J83").Value = 35615
$ws.Range("L83").Value = 106845
$ws.Range("N83").Value = -116829
$ws.Range("H103").Value = 48557.332
$ws.Range("J103").Value = 48557.332
$ws.Range("L103").Value = 48557.332
$ws.Range("N103").Value = -50901.332
$ws.Range("H105").Value = 48606.332
$ws.Range("J105").Value = 48606.332
$ws.Range("L105").Value = 48606.332
$ws.Range("N105").Value = -55594.332
$ws.Range("H126").Value = 3677433
$ws.Range("I126").Value = 5883273
$ws.Range("J126").Value = 1033.3334
$ws.Range("K126").Value = 17649819
$ws.Range("L126").Value = 3100.0002
$ws.Range("M126").Value = -17647349
$ws.Range("N126").Value = -8040.0002
$ws.Range("H132").Value = 1523.742
$ws.Range("I132").Value = 921.86664
$ws.Range("J132").Value = 2088
$ws.Range("K132").Value = 2765.59992
$ws.Range("L132").Value = 6264
$ws.Range("M132").Value = -235.5999199999997
$ws.Range("N132").Value = -11324

